$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FactorID")

# Insert two new rows before the existing row 45 (318000/FundScode) so the
# new FactorID group-1 entries (991004 / 991005) follow 991001-991003.
$ws.Rows.Item(45).Insert()
$ws.Rows.Item(45).Insert()

# New row 45: 收盘价 / 991004 / close / 1 / FundIPOC
$ws.Cells.Item(45, 1).Value = "收盘价"
$ws.Cells.Item(45, 2).Value = 991004
$ws.Cells.Item(45, 3).Value = "close"
$ws.Cells.Item(45, 3).Style = $ws.Cells.Item(44, 3).Style
$ws.Cells.Item(45, 4).Value = 1
$ws.Cells.Item(45, 5).Value = "FundIPOC"

# New row 46: 份额 / 991005 / shares / 1 / FundNegotiableShares3
$ws.Cells.Item(46, 1).Value = "份额"
$ws.Cells.Item(46, 2).Value = 991005
$ws.Cells.Item(46, 3).Value = "shares"
$ws.Cells.Item(46, 3).Style = $ws.Cells.Item(44, 3).Style
$ws.Cells.Item(46, 4).Value = 1
$ws.Cells.Item(46, 5).Value = "FundNegotiableShares3"
